# Auto-generated edit script for cryptos.xlsx update
# Commit: Updated cryptos list on Tue Dec 26 09:51:32 UTC 2023 with GitHub Actions
#
# The sheet stores the Price/Volume columns as literal text (inlineStr) even
# though many of the values look numeric (e.g. "1.00", "113.21"). Excel's
# COM automation auto-detects numeric-looking strings and converts them to
# real numbers when assigned via .Value on a "General" formatted cell, so
# for any replacement text that would parse as a number we first force the
# cell's NumberFormat to Text ("@") to keep it a literal string, matching
# the source data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '42.601.63'
$ws.Range('E2').Value = '  -1.27%  '
$ws.Range('D3').Value = '2.226.26'
$ws.Range('E3').Value = '  -2.09%  '
$ws.Range('D4').NumberFormat = "@"
$ws.Range('D4').Value = '1.00'
$ws.Range('E4').Value = '  +0.22%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '113.21'
$ws.Range('E5').Value = '  +1.69%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '273.75'
$ws.Range('E6').Value = '  +3.75%  '
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '0.625'
$ws.Range('E7').Value = '  -2.70%  '
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '1.01'
$ws.Range('E8').Value = '  +0.26%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.605'
$ws.Range('E9').Value = '  -0.15%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '46.19'
$ws.Range('E10').Value = '  -0.47%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.0924'
$ws.Range('E11').Value = '  -0.90%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '8.97'
$ws.Range('E12').Value = '  -1.36%  '
$ws.Range('E13').Value = '  -3.50%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '15.17'
$ws.Range('E14').Value = '  -1.20%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '0.867'
$ws.Range('E15').Value = '  +0.86%  '
$ws.Range('D16').Value = '2.568.07'
$ws.Range('E16').Value = '  -1.65%  '
$ws.Range('D17').Value = '2.237.73'
$ws.Range('E17').Value = '  -1.03%  '
$ws.Range('D18').Value = '42.670.82'
$ws.Range('E18').Value = '  -1.07%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '0.0000106'
$ws.Range('E19').Value = '  -1.11%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '6.73'
$ws.Range('E20').Value = '  +0.33%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '71.88'
$ws.Range('E21').Value = '  +0.06%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '2.32'
$ws.Range('E22').Value = '  -4.75%  '
$ws.Range('B23').Value = 'BitcoinCash'
$ws.Range('C23').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '230.82'
$ws.Range('E23').Value = '  -1.17%  '
$ws.Range('B24').Value = 'PancakeSwap'
$ws.Range('C24').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '2.95'
$ws.Range('E24').Value = '  +4.07%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '9.17'
$ws.Range('E25').Value = '  -1.82%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '11.89'
$ws.Range('E26').Value = '  +5.44%  '
$ws.Range('E27').Value = '  -1.69%  '
$ws.Range('E28').Value = '  -0.80%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '39.92'
$ws.Range('E29').Value = '  -3.66%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '2.23'
$ws.Range('E30').Value = '  -0.39%  '
$ws.Range('E31').Value = '  -2.25%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '173.04'
$ws.Range('E32').Value = '  -0.07%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '20.96'
$ws.Range('E33').Value = '  -1.91%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '0.0885'
$ws.Range('E34').Value = '  -1.07%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '5.53'
$ws.Range('E35').Value = '  -1.46%  '
$ws.Range('E36').Value = '  +12.30%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '0.127'
$ws.Range('E37').Value = '  -2.40%  '
$ws.Range('B38').Value = 'VeChain'
$ws.Range('C38').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '0.0369'
$ws.Range('E38').Value = '  -0.34%  '
$ws.Range('B39').Value = 'RenderToken'
$ws.Range('C39').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '4.61'
$ws.Range('E39').Value = '  -0.03%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '0.105'
$ws.Range('E40').Value = '  +0.92%  '
$ws.Range('E41').Value = '  -1.23%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '70.76'
$ws.Range('E42').Value = '  -6.50%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '13.19'
$ws.Range('E43').Value = '  -7.65%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '0.231'
$ws.Range('E44').Value = '  -2.04%  '
$ws.Range('E45').Value = '  +0.15%  '
$ws.Range('E46').Value = '  -2.59%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '5.57'
$ws.Range('E47').Value = '  -8.22%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '1.26'
$ws.Range('E48').Value = '  +0.64%  '
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '8.42'
$ws.Range('E49').Value = '  -0.78%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '0.0985'
$ws.Range('E50').Value = '  -0.95%  '
$ws.Range('B51').Value = 'Aave'
$ws.Range('C51').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '99.56'
$ws.Range('E51').Value = '  -1.04%  '
